$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.5606716533
$ws.Range("C2").Value = -224.73150253
$ws.Range("D2").Value = -225.29217418
$ws.Range("E2").Value = -224.5064316978

$ws.Range("B3").Value = -0.5692607282
$ws.Range("C3").Value = -224.66826609
$ws.Range("D3").Value = -225.23752682
$ws.Range("E3").Value = -224.5064316978

$ws.Range("B4").Value = -0.5704422631
$ws.Range("C4").Value = -224.63991225
$ws.Range("D4").Value = -225.21035452
$ws.Range("E4").Value = -224.5064316978
